$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.026.92"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.814.68"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'310.23"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.5000"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").Value = "'0.3924"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "'0.09901"
$ws.Range("E9").Value = "  +26.46%  "
$ws.Range("D10").Value = "'1.103"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'40.93"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "'6.421"
$ws.Range("E12").Value = "  +3.33%  "
$ws.Range("D13").Value = "'20.51"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "'0.9999"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "1.805.45"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "'7.271"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "'0.00001143"
$ws.Range("E17").Value = "  +6.34%  "
$ws.Range("D18").Value = "'92.37"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "'0.06646"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'0.9995"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'17.18"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").Value = "'5.934"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "28.073.94"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "'11.08"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "'2.261"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "'158.50"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").Value = "2.020.06"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").Value = "'20.59"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("D29").Value = "'2.396"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").Value = "'127.30"
$ws.Range("E30").Value = "  +3.01%  "
$ws.Range("D31").Value = "'0.1064"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "'1.034"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "'5.573"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "'3.599"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "'0.06738"
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("D38").Value = "'0.2141"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'4.945"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'11.30"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").Value = "'0.6183"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").Value = "'1.174"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").Value = "'0.9992"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'13.16"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'0.5900"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "'3.691"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'1.280"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").Value = "'124.07"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'1.928"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("D50").Value = "'1.177"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "'0.06782"
$ws.Range("E51").Value = "  -0.46%  "

# Row 36/37: FraxShare and VeChain swap positions with updated price/volume
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'8.910"
$ws.Range("E36").Value = "  +1.82%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02331"
$ws.Range("E37").Value = "  +1.15%  "
